$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bill of materials")

# Update quantity for "Pro Micro 3.3V" (row 2) from 2 to 1; formula in D2 recalcs automatically.
$ws.Range("B2").Value = 1

# Add the explanatory note in row 13 (A13), with rich text formatting:
# "Note:" bold+underline, followed by regular explanation text.
$noteCell = $ws.Range("A13")
$noteCell.Value = "Note: I (Matthew) already have a 5v Arduino Pro Micro that we can use for the ground instead of buying a second 3.3v version (hence the quantity being only 1). The reason for the 3.3v version being necessary for the rocket is that LiPo batteries provide just 3.7v and all the sensors onboard and the radio module are compatible with 3.3v logic."

# Apply the built-in "Note" cell style.
$noteCell.Style = "Note"

# Alignment: left/top with wrap text.
$noteCell.HorizontalAlignment = -4131  # xlLeft
$noteCell.VerticalAlignment = -4160    # xlTop
$noteCell.WrapText = $true

# Make the "Note:" portion of the text bold and underlined.
$rt = $noteCell.Characters(1, 5)
$rt.Font.Bold = $true
$rt.Font.Underline = $true
$rt.Font.Size = 11
$rt.Font.Name = "Calibri"

# Touch the remainder run's font too so it carries explicit font info like the rest of the sheet.
$rt2 = $noteCell.Characters(6, 339)
$rt2.Font.Size = 11
$rt2.Font.Name = "Calibri"

# Row height for the note row.
$ws.Rows.Item(13).RowHeight = 151.5

# Sheet view adjustments: scroll back to column A, select D13 (as in the final saved file).
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("D13").Select()

$wb.Save()
